$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename DK1 -> DK, DK1_Central -> DK_Central
$ws.Range("A2").Value = "DK"
$ws.Range("B2").Value = "DK_Central"
$ws.Range("C2").Value = "DK"
$ws.Range("D2").Value = "DK"
$ws.Range("E2").Value = "DK_Central"
$ws.Range("F2").Value = "DK"

# Row 3: clear A3, C3, D3; set B3/E3 to DK_Decentral, F3 to DK
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "DK_Decentral"
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "DK_Decentral"
$ws.Range("F3").Value = "DK"

# Rows 4-7: clear entirely (previously held DK1_SmallDecentral/DK2/etc rows)
$ws.Range("A4:F7").ClearContents()

# Update selection to match target state
$ws.Range("F4:F7").Select()
